$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -72.419
$ws.Range("C2").Value = 18.12
$ws.Range("D2").Value = 30.104
$ws.Range("E2").Value = 5.624
$ws.Range("F2").Value = 5.051
$ws.Range("G2").Value = 1.008
$ws.Range("H2").Value = -74.071
$ws.Range("I2").Value = -18.12
$ws.Range("J2").Value = 30.104
$ws.Range("K2").Value = 5.624
$ws.Range("L2").Value = 3.04
$ws.Range("M2").Value = 0.621
$ws.Range("N2").Value = -23.98
$ws.Range("O2").Value = 6
$ws.Range("P2").Value = 119.619
$ws.Range("Q2").Value = 5.624
$ws.Range("R2").Value = 5.051
$ws.Range("S2").Value = 0.497

$ws.Range("B3").Value = -72.419
$ws.Range("C3").Value = 18.12
$ws.Range("D3").Value = 30.104
$ws.Range("E3").Value = 5.624
$ws.Range("F3").Value = 5.051
$ws.Range("G3").Value = 1.008
$ws.Range("H3").Value = -74.071
$ws.Range("I3").Value = -18.12
$ws.Range("J3").Value = 30.104
$ws.Range("K3").Value = 5.624
$ws.Range("L3").Value = 3.04
$ws.Range("M3").Value = 0.621
$ws.Range("N3").Value = -23.98
$ws.Range("O3").Value = 6
$ws.Range("P3").Value = 119.619
$ws.Range("Q3").Value = 5.624
$ws.Range("R3").Value = 5.051
$ws.Range("S3").Value = 0.497

$ws.Range("B4").Value = -94.486
$ws.Range("C4").Value = 18.12
$ws.Range("D4").Value = 30.008
$ws.Range("E4").Value = 5.624
$ws.Range("F4").Value = 4.831
$ws.Range("G4").Value = 0.761
$ws.Range("H4").Value = -93.742
$ws.Range("I4").Value = -18.12
$ws.Range("J4").Value = 30.008
$ws.Range("K4").Value = 5.624
$ws.Range("L4").Value = 3.228
$ws.Range("M4").Value = 0.527
$ws.Range("N4").Value = -31.287
$ws.Range("O4").Value = 6
$ws.Range("P4").Value = 119.237
$ws.Range("Q4").Value = 5.624
$ws.Range("R4").Value = 4.831
$ws.Range("S4").Value = 0.475

$ws.Range("B5").Value = -94.624
$ws.Range("C5").Value = 18.12
$ws.Range("D5").Value = 30.009
$ws.Range("E5").Value = 5.624
$ws.Range("F5").Value = 4.83
$ws.Range("G5").Value = 0.76
$ws.Range("H5").Value = -95.157
$ws.Range("I5").Value = -18.12
$ws.Range("J5").Value = 30.009
$ws.Range("K5").Value = 5.624
$ws.Range("L5").Value = 3.239
$ws.Range("M5").Value = 0.521
$ws.Range("N5").Value = -31.333
$ws.Range("O5").Value = 6
$ws.Range("P5").Value = 119.241
$ws.Range("Q5").Value = 5.624
$ws.Range("R5").Value = 4.83
$ws.Range("S5").Value = 0.474

$ws.Range("B6").Value = -93.96299999999999
$ws.Range("C6").Value = 18.12
$ws.Range("D6").Value = 30.468
$ws.Range("E6").Value = 5.624
$ws.Range("F6").Value = 4.836
$ws.Range("G6").Value = 0.764
$ws.Range("H6").Value = -94.96599999999999
$ws.Range("I6").Value = -18.12
$ws.Range("J6").Value = 30.468
$ws.Range("K6").Value = 5.624
$ws.Range("L6").Value = 3.237
$ws.Range("M6").Value = 0.521
$ws.Range("N6").Value = -31.113
$ws.Range("O6").Value = 6
$ws.Range("P6").Value = 121.064
$ws.Range("Q6").Value = 5.624
$ws.Range("R6").Value = 4.836
$ws.Range("S6").Value = 0.464

$ws.Range("B7").Value = -93.462
$ws.Range("C7").Value = 18.12
$ws.Range("D7").Value = 53.252
$ws.Range("E7").Value = 5.624
$ws.Range("F7").Value = 4.842
$ws.Range("G7").Value = 0.653
$ws.Range("H7").Value = -94.509
$ws.Range("I7").Value = -18.12
$ws.Range("J7").Value = 53.252
$ws.Range("K7").Value = 5.624
$ws.Range("L7").Value = 3.234
$ws.Range("M7").Value = 0.467
$ws.Range("N7").Value = -30.948
$ws.Range("O7").Value = 6
$ws.Range("P7").Value = 211.597
$ws.Range("Q7").Value = 5.624
$ws.Range("R7").Value = 4.842
$ws.Range("S7").Value = 0.173

$ws.Range("B8").Value = -93.462
$ws.Range("C8").Value = 18.12
$ws.Range("D8").Value = 53.252
$ws.Range("E8").Value = 5.624
$ws.Range("F8").Value = 4.842
$ws.Range("G8").Value = 0.653
$ws.Range("H8").Value = -94.509
$ws.Range("I8").Value = -18.12
$ws.Range("J8").Value = 53.252
$ws.Range("K8").Value = 5.624
$ws.Range("L8").Value = 3.234
$ws.Range("M8").Value = 0.467
$ws.Range("N8").Value = -30.948
$ws.Range("O8").Value = 6
$ws.Range("P8").Value = 211.597
$ws.Range("Q8").Value = 5.624
$ws.Range("R8").Value = 4.842
$ws.Range("S8").Value = 0.173

$ws.Range("B9").Value = -93.96299999999999
$ws.Range("C9").Value = 18.12
$ws.Range("D9").Value = 30.468
$ws.Range("E9").Value = 5.624
$ws.Range("F9").Value = 4.836
$ws.Range("G9").Value = 0.764
$ws.Range("H9").Value = -94.96599999999999
$ws.Range("I9").Value = -18.12
$ws.Range("J9").Value = 30.468
$ws.Range("K9").Value = 5.624
$ws.Range("L9").Value = 3.237
$ws.Range("M9").Value = 0.521
$ws.Range("N9").Value = -31.113
$ws.Range("O9").Value = 6
$ws.Range("P9").Value = 121.064
$ws.Range("Q9").Value = 5.624
$ws.Range("R9").Value = 4.836
$ws.Range("S9").Value = 0.464

$ws.Range("B10").Value = -94.624
$ws.Range("C10").Value = 18.12
$ws.Range("D10").Value = 30.009
$ws.Range("E10").Value = 5.624
$ws.Range("F10").Value = 4.83
$ws.Range("G10").Value = 0.76
$ws.Range("H10").Value = -95.157
$ws.Range("I10").Value = -18.12
$ws.Range("J10").Value = 30.009
$ws.Range("K10").Value = 5.624
$ws.Range("L10").Value = 3.239
$ws.Range("M10").Value = 0.521
$ws.Range("N10").Value = -31.333
$ws.Range("O10").Value = 6
$ws.Range("P10").Value = 119.241
$ws.Range("Q10").Value = 5.624
$ws.Range("R10").Value = 4.83
$ws.Range("S10").Value = 0.474

$ws.Range("B11").Value = -94.486
$ws.Range("C11").Value = 18.12
$ws.Range("D11").Value = 30.008
$ws.Range("E11").Value = 5.624
$ws.Range("F11").Value = 4.831
$ws.Range("G11").Value = 0.761
$ws.Range("H11").Value = -93.742
$ws.Range("I11").Value = -18.12
$ws.Range("J11").Value = 30.008
$ws.Range("K11").Value = 5.624
$ws.Range("L11").Value = 3.228
$ws.Range("M11").Value = 0.527
$ws.Range("N11").Value = -31.287
$ws.Range("O11").Value = 6
$ws.Range("P11").Value = 119.237
$ws.Range("Q11").Value = 5.624
$ws.Range("R11").Value = 4.831
$ws.Range("S11").Value = 0.475
